$excel.DisplayAlerts = $False

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("structure")

# "Implemented Sheet read functionality for insert scripts":
# the personId row's generation technique for insert switches from
# "Random" to the newly supported "Sheet" technique.
$ws.Range("F2").Value = "Sheet"

# Reflect that as the active selection, matching the saved cursor position.
$ws.Range("F2").Select()

# The insert/update/delete helper-sheets were just scratch templates and
# are no longer needed now that the "Sheet" technique is implemented.
$wb.Worksheets("insert").Delete()
$wb.Worksheets("update").Delete()
$wb.Worksheets("delete").Delete()
